$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.02"
$ws.Range("E2").Value = "'-2.17%"
$ws.Range("D3").Value = "'37.88"
$ws.Range("E3").Value = "'-3.97%"
$ws.Range("D4").Value = "'5.049"
$ws.Range("E4").Value = "'-1.62%"
$ws.Range("D5").Value = "'0.07903"
$ws.Range("E5").Value = "'-3.31%"
$ws.Range("D6").Value = "'2.015"
$ws.Range("E6").Value = "'2.46%"
$ws.Range("D7").Value = "'4.386"
$ws.Range("D8").Value = "'8.220"
$ws.Range("E8").Value = "'0.21%"
$ws.Range("E9").Value = "'-2.76%"
$ws.Range("D10").Value = "'0.9259"
$ws.Range("E10").Value = "'-0.17%"
$ws.Range("D11").Value = "'0.1277"
$ws.Range("E11").Value = "'-8.95%"
$ws.Range("D12").Value = "'0.1903"
$ws.Range("E12").Value = "'-3.75%"
$ws.Range("D13").Value = "'0.08748"
$ws.Range("E13").Value = "'-3.57%"
$ws.Range("D14").Value = "'0.03456"
$ws.Range("E14").Value = "'-1.63%"
$ws.Range("D15").Value = "'0.09725"
$ws.Range("E15").Value = "'-1.10%"
$ws.Range("D16").Value = "'0.001394"
$ws.Range("E16").Value = "'-0.17%"
$ws.Range("D17").Value = "'0.006110"
$ws.Range("E17").Value = "'-0.08%"
$ws.Range("D18").Value = "'3.547"
$ws.Range("E18").Value = "'-3.26%"
$ws.Range("D19").Value = "'0.3440"
$ws.Range("E19").Value = "'-0.49%"
$ws.Range("D20").Value = "'0.1300"
$ws.Range("E20").Value = "'-3.32%"
$ws.Range("E21").Value = "'5.54%"
$ws.Range("D22").Value = "'0.2515"
$ws.Range("E22").Value = "'3.67%"
$ws.Range("E23").Value = "'-1.18%"
$ws.Range("D24").Value = "'0.001222"
$ws.Range("E24").Value = "'0.02%"
$ws.Range("D25").Value = "'0.004597"
$ws.Range("E25").Value = "'-3.97%"
$ws.Range("E26").Value = "'176.46%"
$ws.Range("D39").Value = "'0.02248"
$ws.Range("E39").Value = "'3.41%"
$ws.Range("D40").Value = "'0.05010"
$ws.Range("E40").Value = "'-3.72%"
$ws.Range("D41").Value = "'0.007523"
$ws.Range("D42").Value = "'0.009920"
$ws.Range("E42").Value = "'1.25%"
$ws.Range("E43").Value = "'-1.16%"
$ws.Range("D44").Value = "'0.002093"
$ws.Range("E44").Value = "'-0.99%"
$ws.Range("D45").Value = "'0.008541"
$ws.Range("E45").Value = "'-6.42%"
$ws.Range("D46").Value = "'0.00006428"
$ws.Range("E46").Value = "'0.70%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("D48").Value = "'0.003003"
$ws.Range("E48").Value = "'8.56%"
$ws.Range("D49").Value = "'0.001202"
$ws.Range("E49").Value = "'0.21%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.21%"
